$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and volume change (column E) values.
# Some single-token numeric-looking prices are prefixed with a literal apostrophe
# so Excel stores them as exact text (matching the workbook's inlineStr format)
# instead of auto-converting them to floating point numbers, which would lose
# trailing zeros / introduce rounding artifacts.
$ws.Range("D2").Value = '26.961.50'
$ws.Range("E2").Value = '  +1.98%  '
$ws.Range("D3").Value = '1.814.47'
$ws.Range("E3").Value = '  +2.50%  '
$ws.Range("D4").Value = '''1.006'
$ws.Range("E4").Value = '  +0.33%  '
$ws.Range("D5").Value = '''311.27'
$ws.Range("E5").Value = '  +1.54%  '
$ws.Range("E6").Value = '  +0.32%  '
$ws.Range("D7").Value = '''0.4293'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '''0.07227'
$ws.Range("E9").Value = '  +0.38%  '
$ws.Range("D10").Value = '2.165.58'
$ws.Range("E10").Value = '  +21.30%  '
$ws.Range("D11").Value = '''0.8633'
$ws.Range("E11").Value = '  +1.60%  '
$ws.Range("D12").Value = '''21.27'
$ws.Range("E12").Value = '  +4.63%  '
$ws.Range("D13").Value = '''5.405'
$ws.Range("E13").Value = '  +3.22%  '
$ws.Range("D14").Value = '''6.592'
$ws.Range("E14").Value = '  +2.49%  '
$ws.Range("D15").Value = '''0.06936'
$ws.Range("E15").Value = '  +1.01%  '
$ws.Range("D16").Value = '''81.09'
$ws.Range("E16").Value = '  +2.12%  '
$ws.Range("D17").Value = '''1.011'
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("D18").Value = '''0.000008884'
$ws.Range("E18").Value = '  +2.68%  '
$ws.Range("D19").Value = '''1.005'
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("D20").Value = '''15.17'
$ws.Range("E20").Value = '  +1.11%  '
$ws.Range("D21").Value = '27.000.23'
$ws.Range("E21").Value = '  +2.11%  '
$ws.Range("D22").Value = '''5.167'
$ws.Range("E22").Value = '  +1.22%  '
$ws.Range("D23").Value = '2.405.60'
$ws.Range("E23").Value = '  +20.07%  '
$ws.Range("D24").Value = '''11.00'
$ws.Range("E24").Value = '  -2.37%  '
$ws.Range("D25").Value = '''153.81'
$ws.Range("E25").Value = '  +1.08%  '
$ws.Range("D26").Value = '''1.881'
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("D27").Value = '''18.28'
$ws.Range("E27").Value = '  +0.88%  '
$ws.Range("D28").Value = '''5.216'
$ws.Range("E28").Value = '  +2.45%  '
$ws.Range("D29").Value = '''1.896'
$ws.Range("E29").Value = '  +9.93%  '
$ws.Range("D30").Value = '''114.54'
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").Value = '''0.08932'
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("D32").Value = '''1.188'
$ws.Range("E32").Value = '  +6.39%  '
$ws.Range("D33").Value = '''0.7450'
$ws.Range("E33").Value = '  +2.78%  '
$ws.Range("D34").Value = '''4.409'
$ws.Range("E34").Value = '  +2.01%  '
$ws.Range("E35").Value = '  +2.26%  '
$ws.Range("D36").Value = '''1.005'
$ws.Range("E36").Value = '  +0.29%  '
$ws.Range("D37").Value = '''1.130'
$ws.Range("E37").Value = '  +4.82%  '
$ws.Range("D38").Value = '''0.05204'
$ws.Range("D39").Value = '''0.01919'
$ws.Range("E39").Value = '  +1.36%  '
$ws.Range("D40").Value = '''0.5093'
$ws.Range("E40").Value = '  +3.32%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '''2.743'
$ws.Range("E41").Value = '  +7.87%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '''0.1651'
$ws.Range("E42").Value = '  +3.00%  '
$ws.Range("D43").Value = '''6.460'
$ws.Range("E43").Value = '  +3.46%  '
$ws.Range("D44").Value = '''8.302'
$ws.Range("E44").Value = '  +3.43%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''10.43'
$ws.Range("E45").Value = '  +2.52%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '''106.45'
$ws.Range("E46").Value = '  +1.43%  '
$ws.Range("E47").Value = '  +0.39%  '
$ws.Range("D48").Value = '''0.4560'
$ws.Range("E48").Value = '  +1.74%  '
$ws.Range("D49").Value = '''1.642'
$ws.Range("E49").Value = '  +3.42%  '
$ws.Range("D50").Value = '''0.06206'
$ws.Range("E51").Value = '  +5.05%  '
